$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Text correction: clarify that the "user manual" is a page within the site.
$ws.Range("D3").Value = "Added a user manual page to help with people who may find the website complicated. Added student names, numbers and courses. "

# Row heights settle to a slightly shorter wrap after the edit.
$ws.Rows.Item(2).RowHeight = 28
$ws.Rows.Item(3).RowHeight = 42

# Selection left on D4 when the file was saved.
$ws.Range("D4").Select()
